$d = $word.ActiveDocument

# Title (Heading1) and the bold "Play ... for Free" line near the end share identical text
$d.Content.Find.Execute("Play Lord Merlin and the Lady of the Lake for Free", $true, $false, $false, $false, $false, $true, 1, $false, "Play Lord Merlin and the Lady of the Lake Free - Review", 2)

# "What we like" bullet list
$d.Content.Find.Execute("Stunning graphics and sound effects", $true, $false, $false, $false, $false, $true, 1, $false, "Stunning and immersive graphics", 2)
$d.Content.Find.Execute("Simple yet unique gameplay", $true, $false, $false, $false, $false, $true, 1, $false, "Exceptional soundtrack", 2)
$d.Content.Find.Execute("Free spins feature with 2 variations", $true, $false, $false, $false, $false, $true, 1, $false, "Satisfying animations", 2)
$d.Content.Find.Execute("High volatility for significant wins", $true, $false, $false, $false, $false, $true, 1, $false, "Unique gameplay features", 2)

# "What we don't like" bullet list
$d.Content.Find.Execute("No progressive jackpot", $true, $false, $false, $false, $false, $true, 1, $false, "Limited free spin variations", 2)

# Closing italic summary paragraph
$d.Content.Find.Execute("Read our review of the Lord Merlin and the Lady of the Lake slot game by Play'n Go and play for free. Featuring stunning graphics, immersive sound effects, and a free spins feature.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Lord Merlin and the Lady of the Lake and play for free. Discover stunning graphics, immersive gameplay, and unique features.", 2)
